$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update REMARKS text for the affected dates (column P)
$ws.Range("P5").Value = "~OB Others|Pentstar's Barter CX GO LIVE| R"
$ws.Range("P6").Value = "~OT ~ = Pentstar hard drive configuration and synchronization ~OB Others|Pentstar's Barter CX GO LIVE| R"
$ws.Range("P7").Value = """ ~OT ~ = Barter CX Pentstar Rollout: Robinsons Ermita, SM Mall of Asia ~OB Others|Pentstar's Barter CX GO LIVE| R """
$ws.Range("P10").Value = "~OT ~ = Barter CX Pentstar Rollout - SM Davao and Abreeza Mall ~OB Others|Pentstar's Barter CX GO LIVE| R"
$ws.Range("P11").Value = "~OT ~ = Barter CX Pentstar Rollout - Veranza Mall Gensan ~OB Others|Pentstar's Barter CX GO LIVE| R"
$ws.Range("P12").Value = "~OT ~ = Barter CX Pentstar Rollout - Ayala Center Cebu ~OB Others|Pentstar's Barter CX GO LIVE| R"

# Update overtime hours (column G)
$ws.Range("G7").Value = 5.5
$ws.Range("G11").Value = 5.5

# Update sick leave hours (column I)
$ws.Range("I18").Value = 0.5
$ws.Range("I23").Value = 1.0
